# Auto update Excel log
# Appends new sensor-log rows to the "PIR" sheet (292-304) and the
# "Humidity" sheet (203-211), matching the source system's export format.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PIR sheet: add rows 292..304
# ---------------------------------------------------------------------
$pir = $wb.Worksheets.Item("PIR")

$pirTimes = @(
    "17:31:18",
    "17:31:19",
    "17:31:24",
    "17:31:30",
    "17:31:35",
    "17:31:40",
    "17:31:45",
    "17:31:50",
    "17:31:55",
    "17:32:00",
    "17:32:05",
    "17:32:10",
    "17:32:15"
)

$pirStartRow = 292

# Pre-format column A as Text so date-looking strings are not
# auto-converted into date serial numbers.
$pirEndRow = $pirStartRow + $pirTimes.Length - 1
$pir.Range("A$pirStartRow`:A$pirEndRow").NumberFormat = "@"

for ($i = 0; $i -lt $pirTimes.Length; $i++) {
    $row = $pirStartRow + $i
    $pir.Cells.Item($row, 1).Value = "2026-01-30"
    $pir.Cells.Item($row, 2).Value = $pirTimes[$i]
    $pir.Cells.Item($row, 3).Value = "17:00"
    $pir.Cells.Item($row, 4).Value = "Bathroom"
    $pir.Cells.Item($row, 5).Value = "No Motion"
    $pir.Cells.Item($row, 6).Value = "Inactive"
}

# ---------------------------------------------------------------------
# Humidity sheet: add rows 203..211
# ---------------------------------------------------------------------
$hum = $wb.Worksheets.Item("Humidity")

$humTimes = @(
    "17:31:19",
    "17:31:30",
    "17:31:35",
    "17:31:40",
    "17:31:50",
    "17:31:55",
    "17:32:00",
    "17:32:10",
    "17:32:16"
)
$humValues = @(
    "87.0%",
    "87.0%",
    "87.0%",
    "87.0%",
    "87.1%",
    "87.1%",
    "87.1%",
    "87.0%",
    "87.1%"
)

$humStartRow = 203
$humEndRow = $humStartRow + $humTimes.Length - 1

# Pre-format columns A (date) and E (percentage) as Text so the
# values are stored literally instead of being parsed into a date
# serial / numeric percentage.
$hum.Range("A$humStartRow`:A$humEndRow").NumberFormat = "@"
$hum.Range("E$humStartRow`:E$humEndRow").NumberFormat = "@"

for ($i = 0; $i -lt $humTimes.Length; $i++) {
    $row = $humStartRow + $i
    $hum.Cells.Item($row, 1).Value = "2026-01-30"
    $hum.Cells.Item($row, 2).Value = $humTimes[$i]
    $hum.Cells.Item($row, 3).Value = "17:00"
    $hum.Cells.Item($row, 4).Value = "Bathroom"
    $hum.Cells.Item($row, 5).Value = $humValues[$i]
    $hum.Cells.Item($row, 6).Value = "Active"
}
